$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed values for rows 2-51: (row, C = IPC PO, D = DELTA, E = DELTA^2)
$rows = @(
    ,@(2, 28.61942153242016, -0.3605784675798382, 0.1300168312822245)
    ,@(3, 28.28017903452562, -0.8698209654743749, 0.7565885119787737)
    ,@(4, 28.88646415724472, -0.4635358427552845, 0.2148654775188518)
    ,@(5, 28.79124348772752, -0.5787565122724807, 0.3349591004978061)
    ,@(6, 28.69696473790746, -0.8430352620925348, 0.7107084531314288)
    ,@(7, 29.86126388352603, 0.3112638835260277, 0.09688520518770456)
    ,@(8, 29.59043104121293, -0.1595689587870659, 0.02546225260838833)
    ,@(9, 30.47360527894018, 0.6336052789401769, 0.4014556495008594)
    ,@(10, 30.35429894197885, 0.5442989419788482, 0.2962613382392936)
    ,@(11, 30.19960459309228, 0.2796045930922801, 0.07817872847829956)
    ,@(12, 29.54863326317109, -0.4313667368289096, 0.1860772616424217)
    ,@(13, 29.85653219908905, -0.1834678009109503, 0.03366043397110009)
    ,@(14, 29.40461724422588, -0.8053827557741222, 0.6486413832983194)
    ,@(15, 29.82918319365324, -0.3908168063467627, 0.1527377761230831)
    ,@(16, 29.64636299647648, -0.7336370035235191, 0.538223252938968)
    ,@(17, 30.33773402865503, -0.1022659713449734, 0.01045832889513092)
    ,@(18, 30.36049336566984, -0.1195066343301612, 0.01428183564892286)
    ,@(19, 30.65505704903011, -0.03494295096988864, 0.001221009822484041)
    ,@(20, 30.32258370656307, -0.4274162934369308, 0.1826846878953645)
    ,@(21, 30.73581033191036, -0.2041896680896436, 0.0416934205545588)
    ,@(22, 31.06682419096773, 0.1168241909677299, 0.01364789159526462)
    ,@(23, 31.71203661379989, 0.692036613799889, 0.4789146748396167)
    ,@(24, 32.12226057527928, 1.002260575279283, 1.00452626075916)
    ,@(25, 32.26900250503424, 0.989002505034243, 0.9781259549640079)
    ,@(26, 32.11758947928325, 0.7375894792832476, 0.5440382399493323)
    ,@(27, 32.53827326168924, 0.9582732616892429, 0.9182876440685401)
    ,@(28, 32.52980955291806, 0.8798095529180614, 0.774064849405879)
    ,@(29, 33.277970159239, 1.397970159238998, 1.954320566122709)
    ,@(30, 32.85379537904976, 0.5737953790497556, 0.3292411370188528)
    ,@(31, 32.96916922584683, 0.5191692258468237, 0.2695366850663902)
    ,@(32, 33.13367271978353, 0.2836727197835245, 0.080470211949382)
    ,@(33, 33.22534331497836, 0.3253433149783618, 0.1058482726011095)
    ,@(34, 33.21840007594129, 0.1184000759412882, 0.01401857798290281)
    ,@(35, 34.21330507641433, 0.8133050764143306, 0.6614651473213201)
    ,@(36, 33.99777500965789, 0.2977750096578902, 0.08866995637675658)
    ,@(37, 34.04644958257992, -0.05355041742008382, 0.002867647205865216)
    ,@(38, 34.57601968097843, 0.1760196809784276, 0.03098292809174742)
    ,@(39, 34.85083618170182, -0.04916381829818306, 0.002417081029656759)
    ,@(40, 35.9687543924692, 0.6687543924692037, 0.4472324374468536)
    ,@(41, 36.10000759518975, 0.4000075951897486, 0.1600060762094858)
    ,@(42, 35.91907778380174, -0.3809222161982575, 0.145101734793392)
    ,@(43, 36.40238053655096, -0.3976194634490327, 0.1581012377134966)
    ,@(44, 36.34026366487587, -0.959736335124127, 0.9210938329574907)
    ,@(45, 37.45542092048556, -0.4445790795144404, 0.1976505579419071)
    ,@(46, 38.17963449233827, -0.3203655076617338, 0.1026340584993604)
    ,@(47, 38.58431861107851, -0.3156813889214902, 0.09965473931140117)
    ,@(48, 38.90683823342697, -0.4931617665730244, 0.2432085280094262)
    ,@(49, 39.41233102678765, -0.4876689732123509, 0.2378210274339886)
    ,@(50, 39.24870500583931, -0.8512949941606962, 0.7247031670830598)
    ,@(51, 39.33456577636483, -1.265434223635175, 1.601323774347159)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}

# Totals row (52) and MSE row (53)
$ws.Cells.Item(52, 3).Value = -0.008685308628653132
$ws.Cells.Item(52, 5).Value = 18.1450358373095
$ws.Cells.Item(53, 5).Value = 0.36290071674619
